# Trade #10 closed at 2026-02-16 22:52:57 - base_strategy UP +0.000%
# Adds a new trade row (row 11) to both the "All Trades" sheet and the
# per-strategy "base_strategy" sheet - the two sheets kept in sync in this
# workbook.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 11

    $ws.Cells.Item($row, 1).Value  = 10
    # Column B ("Date") looks like a date ("2026-02-16") and Excel would
    # normally auto-convert a literal like that into a date serial number.
    # Write it as a formula first, then paste-special as a value, so it
    # lands as plain text (matching how the rest of the sheet stores
    # dates/times as text) without leaving a formula behind or touching
    # the cell's number format/style.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.Formula = '="2026-02-16"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($row, 3).Value  = "22:52:57"
    $ws.Cells.Item($row, 4).Value  = "base_strategy"
    $ws.Cells.Item($row, 5).Value  = "UP"
    $ws.Cells.Item($row, 6).Value  = 49.999998
    # Column G ("Exit Price") stays blank - trade is still OPEN.
    $ws.Cells.Item($row, 8).Value  = "OPEN"
    $ws.Cells.Item($row, 9).Value  = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    # Column P ("Exit Reason") stays blank - trade is still OPEN.
    $ws.Cells.Item($row, 17).Value = 0
}
